$d = $word.ActiveDocument

# --- Paragraph 1: "BACKGROUND OF THE STUDY" -> centered, bold heading ---
$p1 = $d.Paragraphs(1)
$p1.Alignment = 1          # wdAlignParagraphCenter
$p1.Range.Font.Bold = $true

# --- Paragraph 2: body text -> first-line indent of 720 twips (0.5in / 36pt) ---
$p2 = $d.Paragraphs(2)
$p2.FirstLineIndent = 36   # points; 36pt = 720 twips

# --- Replace the trailing lone-space run (after the _GoBack bookmark) with
#     the new concluding sentence. We scope the Find to a narrow range that
#     covers only that final run's single space so we don't touch any of
#     the other spaces throughout the paragraph. ---
$paraEnd = $p2.Range.End
$tail = $d.Range($paraEnd - 2, $paraEnd - 1)
$tail.Find.Execute(" ", $false, $false, $false, $false, $false, $true, 1, $false, "Therefore, it is a good choice to create a report management to your business or company sales to make it more detailed and informative. ", 2)
